$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# [week8 day4] cold, uncomfortable, want to cry
# New diary entry for 2018-03-22 (serial 43181): add the missing
# row of data right after the last recorded day (2018-03-21).
$ws.Range("A19").Value = 43181
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 1.5
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = "感冒 "
